$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize data columns B (2) through AH (34) to their new (wider) widths.
# The stored OOXML <col width> ends up as ColumnWidth + 5/6, so each
# target width below has 5/6 pre-subtracted to land exactly on the
# integer widths required (8, 7, 9, etc.).
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 6.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(14).ColumnWidth = 6.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(18).ColumnWidth = 6.166666666666667
$ws.Columns.Item(19).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(25).ColumnWidth = 6.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(31).ColumnWidth = 6.166666666666667
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667
$ws.Columns.Item(33).ColumnWidth = 6.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# Replace the sample data in rows 2-5 with the new (second) dataset.
$ws.Range("A2").Value = 45092.50694444445
$ws.Range("B2").Value = 20.178
$ws.Range("C2").Value = 13.652
$ws.Range("D2").Value = 4.068
$ws.Range("E2").Value = 42.752
$ws.Range("F2").Value = 34.691
$ws.Range("G2").Value = 15.879
$ws.Range("H2").Value = 51.007
$ws.Range("I2").Value = 24.432
$ws.Range("J2").Value = 10.251
$ws.Range("K2").Value = 15.607
$ws.Range("L2").Value = 16.868
$ws.Range("M2").Value = 17.597
$ws.Range("N2").Value = 5.069
$ws.Range("O2").Value = 15.79
$ws.Range("P2").Value = 22.094
$ws.Range("Q2").Value = 13.41
$ws.Range("R2").Value = 3.498
$ws.Range("S2").Value = 2.451
$ws.Range("T2").Value = 232.996
$ws.Range("U2").Value = 43.923
$ws.Range("V2").Value = 14.575
$ws.Range("W2").Value = 29.082
$ws.Range("X2").Value = 14.996
$ws.Range("Y2").Value = 3.13
$ws.Range("Z2").Value = 25.159
$ws.Range("AA2").Value = 12.874
$ws.Range("AB2").Value = 11.654
$ws.Range("AC2").Value = 13.651
$ws.Range("AD2").Value = 17.331
$ws.Range("AE2").Value = 3.457
$ws.Range("AF2").Value = 45.214
$ws.Range("AG2").Value = 8.105
$ws.Range("AH2").Value = 18.222
$ws.Range("A3").Value = 45092.51388888889
$ws.Range("B3").Value = 24.021
$ws.Range("C3").Value = 17.339
$ws.Range("D3").Value = 2.064
$ws.Range("E3").Value = 51.919
$ws.Range("F3").Value = 42.489
$ws.Range("G3").Value = 18.903
$ws.Range("H3").Value = 72.117
$ws.Range("I3").Value = 29.086
$ws.Range("J3").Value = 12.736
$ws.Range("K3").Value = 18.975
$ws.Range("L3").Value = 20.81
$ws.Range("M3").Value = 21.827
$ws.Range("N3").Value = 6.038
$ws.Range("O3").Value = 18.798
$ws.Range("P3").Value = 26.631
$ws.Range("Q3").Value = 15.959
$ws.Range("R3").Value = 1.631
$ws.Range("S3").Value = 1.31
$ws.Range("T3").Value = 278.82
$ws.Range("U3").Value = 52.574
$ws.Range("V3").Value = 17.351
$ws.Range("W3").Value = 35.174
$ws.Range("X3").Value = 18.451
$ws.Range("Y3").Value = 3.055
$ws.Range("Z3").Value = 34.987
$ws.Range("AA3").Value = 15.326
$ws.Range("AB3").Value = 13.72
$ws.Range("AC3").Value = 16.102
$ws.Range("AD3").Value = 21.639
$ws.Range("AE3").Value = 1.266
$ws.Range("AF3").Value = 65.511
$ws.Range("AG3").Value = 9.738
$ws.Range("AH3").Value = 21.692
$ws.Range("A4").Value = 45092.52083333334
$ws.Range("B4").Value = 11.05
$ws.Range("C4").Value = 7.84
$ws.Range("D4").Value = 1.141
$ws.Range("E4").Value = 23.862
$ws.Range("F4").Value = 19.38
$ws.Range("G4").Value = 8.696
$ws.Range("H4").Value = 37.523
$ws.Range("I4").Value = 13.38
$ws.Range("J4").Value = 5.809
$ws.Range("K4").Value = 8.559
$ws.Range("L4").Value = 9.591
$ws.Range("M4").Value = 10.003
$ws.Range("N4").Value = 2.78
$ws.Range("O4").Value = 8.647
$ws.Range("P4").Value = 12.214
$ws.Range("Q4").Value = 7.477
$ws.Range("R4").Value = 1.012
$ws.Range("S4").Value = 0.666
$ws.Range("T4").Value = 124.313
$ws.Range("U4").Value = 24.318
$ws.Range("V4").Value = 7.982
$ws.Range("W4").Value = 16.147
$ws.Range("X4").Value = 8.471
$ws.Range("Y4").Value = 1.542
$ws.Range("Z4").Value = 17.461
$ws.Range("AA4").Value = 7.05
$ws.Range("AB4").Value = 6.383
$ws.Range("AC4").Value = 7.479
$ws.Range("AD4").Value = 9.928
$ws.Range("AE4").Value = 0.773
$ws.Range("AF4").Value = 34.204
$ws.Range("AG4").Value = 4.425
$ws.Range("AH4").Value = 9.979
$ws.Range("A5").Value = 45092.52777777778
$ws.Range("B5").Value = 21.62
$ws.Range("C5").Value = 15.9
$ws.Range("D5").Value = 1.28
$ws.Range("E5").Value = 46.92
$ws.Range("F5").Value = 38.5
$ws.Range("G5").Value = 17.01
$ws.Range("H5").Value = 64.26
$ws.Range("I5").Value = 26.18
$ws.Range("J5").Value = 11.6
$ws.Range("K5").Value = 17.22
$ws.Range("L5").Value = 18.85
$ws.Range("M5").Value = 19.83
$ws.Range("N5").Value = 5.43
$ws.Range("O5").Value = 16.92
$ws.Range("P5").Value = 24.05
$ws.Range("Q5").Value = 14.27
$ws.Range("R5").Value = 0.86
$ws.Range("S5").Value = 0.87
$ws.Range("T5").Value = 250.2
$ws.Range("U5").Value = 47.27
$ws.Range("V5").Value = 15.62
$ws.Range("W5").Value = 31.77
$ws.Range("X5").Value = 16.72
$ws.Range("Y5").Value = 2.53
$ws.Range("Z5").Value = 31.32
$ws.Range("AA5").Value = 13.79
$ws.Range("AB5").Value = 12.25
$ws.Range("AC5").Value = 14.4
$ws.Range("AD5").Value = 19.71
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 58.19
$ws.Range("AG5").Value = 8.8
$ws.Range("AH5").Value = 19.52

# Row 6 (the old final data row) is no longer part of the dataset.
$ws.Rows.Item(6).Delete()
